$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 451, shifting existing rows 451:561 down to 452:562
$ws.Rows.Item(451).Insert()

# Populate the newly inserted row 451 with the new weekly price record
$ws.Cells.Item(451, 1).Value = 5
$ws.Cells.Item(451, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(451, 3).Value = "Maule"
$ws.Cells.Item(451, 4).Value = 44754
$ws.Cells.Item(451, 5).Value = 7
$ws.Cells.Item(451, 6).Value = "Fruta"
$ws.Cells.Item(451, 7).Value = 100102
$ws.Cells.Item(451, 8).Value = "Cítricos"
$ws.Cells.Item(451, 9).Value = 100102005
$ws.Cells.Item(451, 10).Value = "Naranja"
$ws.Cells.Item(451, 11).Value = "Fukumoto"
$ws.Cells.Item(451, 12).Value = "Primera"
$ws.Cells.Item(451, 13).Value = 600
$ws.Cells.Item(451, 14).Value = 7000
$ws.Cells.Item(451, 15).Value = 7000
$ws.Cells.Item(451, 16).Value = 7000
$ws.Cells.Item(451, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(451, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(451, 19).Value = 467
$ws.Cells.Item(451, 20).Value = 15
